# Trade #88 (MarketMaking strategy-sheet trade #116 / global trade #116)
# closed at 2026-02-17 21:18:28, plus a brand-new open trade (#149) logged
# right after it. Update the rollup sheets (Summary, Strategy Status) and
# the two trade-log sheets (All Trades, MarketMaking) to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.08   # Current Capital
$summary.Range("B4").Value = 0.87      # Total P&L $
$summary.Range("B6").Value = 116       # Total Trades
$summary.Range("B7").Value = 52        # Winning Trades
$summary.Range("B9").Value = 44.83     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.08     # Capital
$status.Range("D5").Value = 83         # Trades
$status.Range("E5").Value = 0.76       # P&L $
$status.Range("F5").Value = 1.08       # P&L %
$status.Range("G5").Value = 45.78      # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - trade #116 (row 117) transitions OPEN -> CLOSED
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(117, 7).Value = 0.9399999999999999   # G - Exit Price
$allTrades.Cells.Item(117, 8).Value = "CLOSED"              # H - Status
$allTrades.Cells.Item(117, 9).Value = 1.0753                # I - P&L %
$allTrades.Cells.Item(117, 10).Value = 0.01                 # J - P&L $
$allTrades.Cells.Item(117, 11).Value = 101.08                # K - Capital After
$allTrades.Cells.Item(117, 12).Value = "early_exit"          # L - Exit Reason
$allTrades.Cells.Item(117, 13).Value = 0.13                  # M - Duration (min)

# New trade #149 appended as row 150
# (the date string is prefixed with an apostrophe + style reset so Excel
#  stores it as literal text "2026-02-17" instead of auto-converting it
#  to a date serial number, matching the other date cells in this column)
$allTrades.Cells.Item(150, 1).Value = 149
$allTrades.Cells.Item(150, 2).Value = "'2026-02-17"
$allTrades.Cells.Item(150, 2).Style = "Normal"
$allTrades.Cells.Item(150, 3).Value = "21:18:22"
$allTrades.Cells.Item(150, 4).Value = "MarketMaking"
$allTrades.Cells.Item(150, 5).Value = "DOWN"
$allTrades.Cells.Item(150, 6).Value = 0.93
$allTrades.Cells.Item(150, 8).Value = "OPEN"
$allTrades.Cells.Item(150, 9).Value = 0
$allTrades.Cells.Item(150, 10).Value = 0
$allTrades.Cells.Item(150, 11).Value = 101.0696151053151
$allTrades.Cells.Item(150, 13).Value = 0
$allTrades.Cells.Item(150, 14).Value = 0
$allTrades.Cells.Item(150, 15).Value = 0
$allTrades.Cells.Item(150, 16).Value = 0.6
$allTrades.Cells.Item(150, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet - strategy-local trade #116 (row 84) OPEN -> CLOSED
# (column layout differs from "All Trades": L/M/N/O/P/Q =
#  Entry Slippage / Exit Slippage / Confidence / Entry Reason /
#  Exit Reason / Duration)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(84, 7).Value = 0.9399999999999999   # G - Exit Price
$mm.Cells.Item(84, 8).Value = "CLOSED"              # H - Status
$mm.Cells.Item(84, 9).Value = 1.0753                # I - P&L %
$mm.Cells.Item(84, 10).Value = 0.01                 # J - P&L $
$mm.Cells.Item(84, 11).Value = 101.08                # K - Capital After
$mm.Cells.Item(84, 16).Value = "early_exit"          # P - Exit Reason
$mm.Cells.Item(84, 17).Value = 0.13                  # Q - Duration (min)

# New trade #149 appended as row 117
$mm.Cells.Item(117, 1).Value = 149
$mm.Cells.Item(117, 2).Value = "'2026-02-17"
$mm.Cells.Item(117, 2).Style = "Normal"
$mm.Cells.Item(117, 3).Value = "21:18:22"
$mm.Cells.Item(117, 4).Value = "MarketMaking"
$mm.Cells.Item(117, 5).Value = "DOWN"
$mm.Cells.Item(117, 6).Value = 0.93
$mm.Cells.Item(117, 8).Value = "OPEN"
$mm.Cells.Item(117, 9).Value = 0
$mm.Cells.Item(117, 10).Value = 0
$mm.Cells.Item(117, 11).Value = 101.0696151053151
$mm.Cells.Item(117, 12).Value = 0
$mm.Cells.Item(117, 13).Value = 0
$mm.Cells.Item(117, 14).Value = 0.6
$mm.Cells.Item(117, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(117, 17).Value = 0
